$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add a new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the bold/centered/bordered header style used on the other sheets
# (copy the exact cell format rather than rebuilding it by hand, so the
# workbook reuses the existing style record instead of creating a new one)
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$wsForecast.Application.CutCopyMode = $false

# Data rows
$dates = @(45431.99999999999, 45501.99999999999, 45515.99999999999, 45585.99999999999, 45592.99999999999, 45599.99999999999, 45606.99999999999, 45613.99999999999, 45620.99999999999, 45627.99999999999, 45634.99999999999, 45641.99999999999)
$lowers = @(1.999999997302766, 1.999999997030882, 1.999999997452852, 1.999999997358232, 1.99999999736689, 1.999999997189573, 1.999999997149367, 1.999999997056642, 1.999999996577559, 1.99999999603486, 1.999999995320062, 1.999999994948777)
$uppers = @(2.000000002668533, 2.000000002520904, 2.000000002661751, 2.000000002629477, 2.000000002510896, 2.000000002565524, 2.000000002905992, 2.00000000303704, 2.000000003063464, 2.000000003633107, 2.000000003957209, 2.000000004259563)

for ($i = 0; $i -lt $dates.Count; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 1).Value = $dates[$i]
    $wsForecast.Cells.Item($row, 2).Value = 2
    $wsForecast.Cells.Item($row, 3).Value = $lowers[$i]
    $wsForecast.Cells.Item($row, 4).Value = $uppers[$i]
}

# Copy the date-format style used for column A on the other sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A13").PasteSpecial(-4122)
$wsForecast.Application.CutCopyMode = $false
